$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.723979
$ws.Range("H2").Value = 23.171937
$ws.Range("I2").Value = 0.471042132528101
$ws.Range("J2").Value = 0.471042132528101
$ws.Range("Q2").Value = 0.05476816042933333
$ws.Range("R2").Value = 0.492913443864
$ws.Range("S2").Value = 0.471042132528101
$ws.Range("T2").Value = 0.471042132528101

# Row 3
$ws.Range("I3").Value = 0.2460132574367717
$ws.Range("J3").Value = 0.2460132574367717
$ws.Range("S3").Value = 0.2460132574367717
$ws.Range("T3").Value = 0.2460132574367717

# Row 4
$ws.Range("I4").Value = 0.2829446100351274
$ws.Range("J4").Value = 0.2829446100351274
$ws.Range("S4").Value = 0.2829446100351274
$ws.Range("T4").Value = 0.2829446100351274
